$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B20").Value = 6187
$ws.Range("C20").Value = 980
$ws.Range("D20").Value = 5581453
$ws.Range("E20").Value = 902.1259091643769
$ws.Range("F20").Value = 6.875107963378824
$ws.Range("G20").Value = 4.033970276008492
$ws.Range("H20").Value = 26.25799284407599
